$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.416.50'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.917.41'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.10'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4816'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4059'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08211'
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.007'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.34'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.916.71'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.073'
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.232'
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.59'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06873'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.59'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.428.83'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.659'
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.77'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.147.93'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.561'
$ws.Range('E26').Value = '  +7.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.78'
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.98'
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.56'
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09636'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.622'
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.553'
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.372'
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02281'
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06095'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.180'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.92'
$ws.Range('E39').Value = '  +6.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.046'
$ws.Range('E40').Value = '  +1.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5949'
$ws.Range('E41').Value = '  +1.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1846'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.282'
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.375'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07608'
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.44'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5578'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.949'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '118.95'
$ws.Range('E49').Value = '  +4.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.429'
$ws.Range('E50').Value = '  +3.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.11'
$ws.Range('E51').Value = '  -0.32%  '
